# MetodosESB.xlsx update: add report tasks, drop the empty Sheet3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Remove the unused "Sheet3"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
[void]$ws3.Delete()

# ---------------------------------------------------------------------
# 2. Tasks sheet updates
# ---------------------------------------------------------------------
$wsTasks = $wb.Worksheets.Item("Tasks")

# ESBServices is now Done
$wsTasks.Range("B20").Value = "Done"

# Rebuild the "Relatorio" block (rows 25-29 -> rows 25-32) with the
# updated report-task breakdown.
$wsTasks.Range("A25").Value = "diagrama dinâmico - WS+ESB"
$wsTasks.Range("C25").Value = "Joao"
$wsTasks.Range("D25").Value = "Incluir alguma explicação"

$wsTasks.Range("A26").Value = "diagrama dinâmico - BPM"
$wsTasks.Range("C26").Value = "Filipe"
$wsTasks.Range("D26").Value = "BPM print "

$wsTasks.Range("A27").Value = "responsibilities (tabela com packages)"
$wsTasks.Range("C27").Value = "Joao "

$wsTasks.Range("A28").Value = "Design decisions - BPM"
$wsTasks.Range("C28").Value = "Filipe"

$wsTasks.Range("A29").Value = "Design decisions - WS"
$wsTasks.Range("C29").Value = "Joao "
$wsTasks.Range("D29").Value = "Only if necessary"

$wsTasks.Range("A30").Value = "Design decisions - Splitter"
$wsTasks.Range("C30").Value = "David"

$wsTasks.Range("A31").Value = "Deployment instructions - projecto"
$wsTasks.Range("C31").Value = "David"

$wsTasks.Range("A32").Value = "Deployment instructions - jbpm"
$wsTasks.Range("C32").Value = "Filipe"

# ---------------------------------------------------------------------
# 3. View state: Tasks becomes the active tab / sheet with B27 selected
#    (Serviços In Out keeps its existing C6 selection untouched so it
#    doesn't steal the active-sheet state back.)
# ---------------------------------------------------------------------
$wsTasks.Activate()
[void]$wsTasks.Range("B27").Select()
